$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (font/style) from the last existing data row (60),
# which mirrors the same 4-column (A:D) layout, onto the two new rows
# before filling in the new content. This keeps column C's CHEBI-id style
# (s="7") and the other columns' default style (s="1") consistent with the
# rest of the sheet, without touching any unused columns/rows.
$ws.Range("A60:D60").Copy() | Out-Null
$ws.Range("A62:D62").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A60:D60").Copy() | Out-Null
$ws.Range("A63:D63").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# New row 62: cardiovascular agent therapy
$ws.Range("C62").Value = "CHEBI:35554"
$ws.Range("D62").Value = "cardiovascular drug"
$ws.Range("B62").Value = "cardiovascular agent therapy"
$ws.Range("A62").Value = "MAXO_0000181"

# New row 63: antiarrhythmic agent therapy
$ws.Range("A63").Value = "MAXO_0000185"
$ws.Range("B63").Value = "antiarrythmic agent tehrapy"
$ws.Range("C63").Value = "CHEBI:38070"
$ws.Range("D63").Value = "anti-arrhythmia drug"

# Match the row heights (17) used by the neighboring rows 59/60.
$ws.Rows.Item(62).RowHeight = 17
$ws.Rows.Item(63).RowHeight = 17

# Update selection to match the post-edit cursor position.
$ws.Range("D64").Select() | Out-Null
